# EU-27 update for InputData elec/BPHC
# - About sheet: add a second citation row for the UK annual report (now
#   tracked separately post-Brexit), bump the source-row citation from
#   row 20 to row 52, and turn the JRC URL into a live hyperlink.
# - JRC_POTEnCIA sheet: relabel EU28 -> EU27 and refresh the underlying
#   "Net capacities installed (MW)" series with EU27 figures.
# - BPHC sheet formulas reference JRC_POTEnCIA and refresh automatically on
#   recalculation.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsJRC   = $wb.Worksheets.Item("JRC_POTEnCIA")
$wsBPHC  = $wb.Worksheets.Item("BPHC")

# ---------------------------------------------------------------------------
# About sheet - new UK citation row, refreshed row reference, live hyperlink
# ---------------------------------------------------------------------------
$wsAbout.Range("B12").Value = "Annual Reports, Power Generation, Central_2018_UK_pg_det_yearly"

# ---------------------------------------------------------------------------
# JRC_POTEnCIA sheet - EU28 -> EU27 relabel + refreshed series
# ---------------------------------------------------------------------------
$wsJRC.Range("A1").Value = "EU27: Net capacities installed (MW)"

# ---------------------------------------------------------------------------
# About sheet - remaining citation updates (row 20 -> row 52)
# ---------------------------------------------------------------------------
$wsAbout.Range("B11").Value = "`"Net Capacities`", row 52"
$wsAbout.Range("B13").Value = "`"Net Capacities`", row 52"

$wsJRC.Range("B2").Value = 38746.120000000003
$wsJRC.Range("C2").Value = 38918.520000000004
$wsJRC.Range("D2").Value = 38980.520000000004
$wsJRC.Range("E2").Value = 39040.020000000004
$wsJRC.Range("F2").Value = 40096.42
$wsJRC.Range("G2").Value = 40851.980000000003
$wsJRC.Range("H2").Value = 41294.58
$wsJRC.Range("I2").Value = 41294.58
$wsJRC.Range("J2").Value = 41572.58
$wsJRC.Range("K2").Value = 41968.58
$wsJRC.Range("L2").Value = 42327.38
$wsJRC.Range("M2").Value = 42567.38
$wsJRC.Range("N2").Value = 42742.48
$wsJRC.Range("O2").Value = 43171.48
$wsJRC.Range("P2").Value = 43524.480000000003
$wsJRC.Range("Q2").Value = 44591.48
$wsJRC.Range("R2").Value = 45216.08
$wsJRC.Range("S2").Value = 45466.080000000002
$wsJRC.Range("T2").Value = 45858.080000000002
$wsJRC.Range("U2").Value = 45858.080000000002
$wsJRC.Range("V2").Value = 45858.080000000002
$wsJRC.Range("W2").Value = 45858.080000000002
$wsJRC.Range("X2").Value = 45715.08
$wsJRC.Range("Y2").Value = 45615.08
$wsJRC.Range("Z2").Value = 45615.08
$wsJRC.Range("AA2").Value = 45615.08
$wsJRC.Range("AB2").Value = 45608.98
$wsJRC.Range("AC2").Value = 45608.98
$wsJRC.Range("AD2").Value = 45608.98
$wsJRC.Range("AE2").Value = 45608.98
$wsJRC.Range("AF2").Value = 45608.98
$wsJRC.Range("AG2").Value = 45608.98
$wsJRC.Range("AH2").Value = 45465.98
$wsJRC.Range("AI2").Value = 45300.98
$wsJRC.Range("AJ2").Value = 45300.98
$wsJRC.Range("AK2").Value = 45300.98
$wsJRC.Range("AL2").Value = 45300.98
$wsJRC.Range("AM2").Value = 45300.98
$wsJRC.Range("AN2").Value = 45300.98
$wsJRC.Range("AO2").Value = 45300.98
$wsJRC.Range("AP2").Value = 45300.98
$wsJRC.Range("AQ2").Value = 45225.98
$wsJRC.Range("AR2").Value = 45150.98
$wsJRC.Range("AS2").Value = 45150.98
$wsJRC.Range("AT2").Value = 45150.98
$wsJRC.Range("AU2").Value = 45084.98
$wsJRC.Range("AV2").Value = 45084.98
$wsJRC.Range("AW2").Value = 45084.98
$wsJRC.Range("AX2").Value = 44885.38
$wsJRC.Range("AY2").Value = 44833.120000000003
$wsJRC.Range("AZ2").Value = 44744.160000000003

# ---------------------------------------------------------------------------
# About sheet - JRC URL becomes a real hyperlink (keep the existing look)
# ---------------------------------------------------------------------------
$jrcUrl = "https://ec.europa.eu/jrc/en/publication/eur-scientific-and-technical-research-reports/potencia-central-scenario-eu-energy-outlook-2050"
$wsAbout.Hyperlinks.Add($wsAbout.Range("B8"), $jrcUrl)
$wsAbout.Range("B8").Style = $wsAbout.Range("B9").Style

# ---------------------------------------------------------------------------
# Restore the on-disk selections / active sheet seen in the authored file
# ---------------------------------------------------------------------------
$wsJRC.Activate()
$wsJRC.Range("E20").Select()

$wsBPHC.Activate()
$wsBPHC.Range("B2").Select()

$wsAbout.Activate()
$wsAbout.Range("B15").Select()
